# ------------------------------------------------------------------
# fix : fix error due to poor connexion
#
# Populates the "debit" sheet (previously empty) with the filtered
# "debit" rows, fixes up the view/selection state of the "debit" and
# "tout" sheets and widens two columns on "tout" so the long labels in
# them are readable.
# ------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$wsDebit = $wb.Worksheets.Item("debit")
$wsTout  = $wb.Worksheets.Item("tout")

# --- "tout" sheet: widen the Libelle (C) and TableHisto (J) columns ---
$wsTout.Columns.Item(3).ColumnWidth = 31.15
$wsTout.Columns.Item(10).ColumnWidth = 27.5

# --- "tout" sheet: scroll/selection ends up on row 874 (whole row) ---
$wsTout.Rows.Item(874).Select()

# --- "debit" sheet: fill in the 22 x 18 data table (1 header row + 21
#     data rows, columns A:R) that was previously completely empty ---
$arr = New-Object 'object[,]' 22,18
$arr[0,0] = "Tag"
$arr[0,1] = "TagStation"
$arr[0,2] = "Libelle"
$arr[0,3] = "MinStation"
$arr[0,4] = "MaxStation"
$arr[0,5] = "MinPC"
$arr[0,6] = "MaxPC"
$arr[0,7] = "Format"
$arr[0,8] = "Unite"
$arr[0,9] = "TableHisto"
$arr[0,10] = "TableBilan"
$arr[0,11] = "Categorie"
$arr[0,12] = "tag_type_domaine"
$arr[0,13] = "tag_domaine"
$arr[0,14] = "tag_sous_domaine"
$arr[0,15] = "tag_bassin"
$arr[0,16] = "tag_sous_bassin"
$arr[0,17] = "tag_sous_bassin_loc"
$arr[1,0] = 1005
$arr[1,1] = 10
$arr[1,2] = "débit total"
$arr[1,3] = 0
$arr[1,4] = 1
$arr[1,5] = 0
$arr[1,6] = 1
$arr[1,7] = 2
$arr[1,8] = $null
$arr[1,9] = "B_AUCFER_DEBIT"
$arr[1,10] = "C_AUCFER_DEBIT"
$arr[1,11] = 16
$arr[1,12] = 0
$arr[1,13] = 0
$arr[1,14] = 0
$arr[1,15] = 0
$arr[1,16] = 0
$arr[1,17] = 0
$arr[2,0] = 1588
$arr[2,1] = 15
$arr[2,2] = "Débit"
$arr[2,3] = 0
$arr[2,4] = 1
$arr[2,5] = 0
$arr[2,6] = 1
$arr[2,7] = 2
$arr[2,8] = "l/s"
$arr[2,9] = "B_SIPHON_DEBIT"
$arr[2,10] = "C_SIPHON_DEBIT"
$arr[2,11] = 16
$arr[2,12] = 0
$arr[2,13] = 0
$arr[2,14] = 0
$arr[2,15] = 0
$arr[2,16] = 0
$arr[2,17] = 0
$arr[3,0] = 1900
$arr[3,1] = 19
$arr[3,2] = "Débit moyen"
$arr[3,3] = 0
$arr[3,4] = 1
$arr[3,5] = 0
$arr[3,6] = 1
$arr[3,7] = 2
$arr[3,8] = "m³/s"
$arr[3,9] = "B_PONT_DE_CRAN_DEBIT"
$arr[3,10] = "C_PONT_DE_CRAN_DEBIT"
$arr[3,11] = 4
$arr[3,12] = 1
$arr[3,13] = 3
$arr[3,14] = 0
$arr[3,15] = 1
$arr[3,16] = 1
$arr[3,17] = 4
$arr[4,0] = 1903
$arr[4,1] = 19
$arr[4,2] = "Débit corde inférieure ( < 03/2009 )"
$arr[4,3] = 0
$arr[4,4] = 1
$arr[4,5] = 0
$arr[4,6] = 1
$arr[4,7] = 2
$arr[4,8] = $null
$arr[4,9] = "B_PONT_DE_CRAN_DEBIT"
$arr[4,10] = "C_PONT_DE_CRAN_DEBIT"
$arr[4,11] = 4
$arr[4,12] = 1
$arr[4,13] = 3
$arr[4,14] = 0
$arr[4,15] = 1
$arr[4,16] = 1
$arr[4,17] = 4
$arr[5,0] = 1904
$arr[5,1] = 19
$arr[5,2] = "Débit corde intermédiaire ( < 03/2009 )"
$arr[5,3] = 0
$arr[5,4] = 1
$arr[5,5] = 0
$arr[5,6] = 1
$arr[5,7] = 2
$arr[5,8] = $null
$arr[5,9] = "B_PONT_DE_CRAN_DEBIT"
$arr[5,10] = "C_PONT_DE_CRAN_DEBIT"
$arr[5,11] = 4
$arr[5,12] = 1
$arr[5,13] = 3
$arr[5,14] = 0
$arr[5,15] = 1
$arr[5,16] = 1
$arr[5,17] = 4
$arr[6,0] = 1905
$arr[6,1] = 19
$arr[6,2] = "Débit corde supérieure ( < 03/2009 )"
$arr[6,3] = 0
$arr[6,4] = 1
$arr[6,5] = 0
$arr[6,6] = 1
$arr[6,7] = 2
$arr[6,8] = $null
$arr[6,9] = "B_PONT_DE_CRAN_DEBIT"
$arr[6,10] = "C_PONT_DE_CRAN_DEBIT"
$arr[6,11] = 4
$arr[6,12] = 1
$arr[6,13] = 3
$arr[6,14] = 0
$arr[6,15] = 1
$arr[6,16] = 1
$arr[6,17] = 4
$arr[7,0] = 2515
$arr[7,1] = 25
$arr[7,2] = "Débit Vilaine estimé"
$arr[7,3] = 0
$arr[7,4] = 1
$arr[7,5] = 0
$arr[7,6] = 1
$arr[7,7] = 2
$arr[7,8] = "m³/s"
$arr[7,9] = "B_BARRAGE_DEBIT"
$arr[7,10] = "C_BARRAGE_DEBIT"
$arr[7,11] = 4
$arr[7,12] = 2
$arr[7,13] = 3
$arr[7,14] = 0
$arr[7,15] = 1
$arr[7,16] = 1
$arr[7,17] = 0
$arr[8,0] = 2523
$arr[8,1] = 25
$arr[8,2] = "Débit passe"
$arr[8,3] = 0
$arr[8,4] = 1
$arr[8,5] = 0
$arr[8,6] = 1
$arr[8,7] = 2
$arr[8,8] = "m³/s"
$arr[8,9] = "B_BARRAGE_DEBIT"
$arr[8,10] = "C_BARRAGE_DEBIT"
$arr[8,11] = 4
$arr[8,12] = 2
$arr[8,13] = 3
$arr[8,14] = 0
$arr[8,15] = 1
$arr[8,16] = 1
$arr[8,17] = 0
$arr[9,0] = 2536
$arr[9,1] = 25
$arr[9,2] = "Volume/jour Vannes"
$arr[9,3] = 0
$arr[9,4] = 1
$arr[9,5] = 0
$arr[9,6] = 1
$arr[9,7] = "."
$arr[9,8] = "m³"
$arr[9,9] = "B_BARRAGE_VOLUME"
$arr[9,10] = "C_BARRAGE_VOLUME"
$arr[9,11] = 3
$arr[9,12] = 2
$arr[9,13] = 2
$arr[9,14] = 0
$arr[9,15] = 1
$arr[9,16] = 1
$arr[9,17] = 0
$arr[10,0] = 2537
$arr[10,1] = 25
$arr[10,2] = "Volume/jour Passe"
$arr[10,3] = 0
$arr[10,4] = 1
$arr[10,5] = 0
$arr[10,6] = 1
$arr[10,7] = 0
$arr[10,8] = "m³"
$arr[10,9] = "B_BARRAGE_VOLUME"
$arr[10,10] = "C_BARRAGE_VOLUME"
$arr[10,11] = 3
$arr[10,12] = 2
$arr[10,13] = 2
$arr[10,14] = 0
$arr[10,15] = 1
$arr[10,16] = 1
$arr[10,17] = 0
$arr[11,0] = 2538
$arr[11,1] = 25
$arr[11,2] = "Volume/jour Ecluse"
$arr[11,3] = 0
$arr[11,4] = 1
$arr[11,5] = 0
$arr[11,6] = 1
$arr[11,7] = 0
$arr[11,8] = "m³"
$arr[11,9] = "B_BARRAGE_VOLUME"
$arr[11,10] = "C_BARRAGE_VOLUME"
$arr[11,11] = 3
$arr[11,12] = 2
$arr[11,13] = 2
$arr[11,14] = 0
$arr[11,15] = 1
$arr[11,16] = 1
$arr[11,17] = 0
$arr[12,0] = 2539
$arr[12,1] = 25
$arr[12,2] = "Volume/jour Siphon"
$arr[12,3] = 0
$arr[12,4] = 1
$arr[12,5] = 0
$arr[12,6] = 1
$arr[12,7] = 0
$arr[12,8] = "m³"
$arr[12,9] = "B_BARRAGE_VOLUME"
$arr[12,10] = "C_BARRAGE_VOLUME"
$arr[12,11] = 3
$arr[12,12] = 2
$arr[12,13] = 2
$arr[12,14] = 0
$arr[12,15] = 1
$arr[12,16] = 1
$arr[12,17] = 0
$arr[13,0] = 2540
$arr[13,1] = 25
$arr[13,2] = "Volume/jour Volets"
$arr[13,3] = 0
$arr[13,4] = 1
$arr[13,5] = 0
$arr[13,6] = 1
$arr[13,7] = -2
$arr[13,8] = "m³"
$arr[13,9] = "B_BARRAGE_VOLUME"
$arr[13,10] = "C_BARRAGE_VOLUME"
$arr[13,11] = 3
$arr[13,12] = 2
$arr[13,13] = 2
$arr[13,14] = 0
$arr[13,15] = 1
$arr[13,16] = 1
$arr[13,17] = 0
$arr[14,0] = 2550
$arr[14,1] = 25
$arr[14,2] = "Totalisateur Volumes évacués vannes"
$arr[14,3] = 0
$arr[14,4] = 1
$arr[14,5] = 0
$arr[14,6] = 1
$arr[14,7] = 0
$arr[14,8] = "m³"
$arr[14,9] = "B_BARRAGE_VOLUME"
$arr[14,10] = "C_BARRAGE_VOLUME"
$arr[14,11] = 3
$arr[14,12] = 2
$arr[14,13] = 2
$arr[14,14] = 0
$arr[14,15] = 1
$arr[14,16] = 1
$arr[14,17] = 0
$arr[15,0] = 2551
$arr[15,1] = 25
$arr[15,2] = "Totalisateur Volumes évacués passe"
$arr[15,3] = 0
$arr[15,4] = 1
$arr[15,5] = 0
$arr[15,6] = 1
$arr[15,7] = 2
$arr[15,8] = "m³"
$arr[15,9] = "B_BARRAGE_VOLUME"
$arr[15,10] = "C_BARRAGE_VOLUME"
$arr[15,11] = 3
$arr[15,12] = 2
$arr[15,13] = 2
$arr[15,14] = 0
$arr[15,15] = 1
$arr[15,16] = 1
$arr[15,17] = 0
$arr[16,0] = 2552
$arr[16,1] = 25
$arr[16,2] = "Totalisateur Volumes évacués siphon"
$arr[16,3] = 0
$arr[16,4] = 1
$arr[16,5] = 0
$arr[16,6] = 1
$arr[16,7] = 0
$arr[16,8] = "m³"
$arr[16,9] = "B_BARRAGE_VOLUME"
$arr[16,10] = "C_BARRAGE_VOLUME"
$arr[16,11] = 3
$arr[16,12] = 2
$arr[16,13] = 2
$arr[16,14] = 0
$arr[16,15] = 1
$arr[16,16] = 1
$arr[16,17] = 0
$arr[17,0] = 2553
$arr[17,1] = 25
$arr[17,2] = "Totalisateur Volumes évacués volet"
$arr[17,3] = 0
$arr[17,4] = 1
$arr[17,5] = 0
$arr[17,6] = 1
$arr[17,7] = 0
$arr[17,8] = "m³"
$arr[17,9] = "B_BARRAGE_VOLUME"
$arr[17,10] = "C_BARRAGE_VOLUME"
$arr[17,11] = 3
$arr[17,12] = 2
$arr[17,13] = 2
$arr[17,14] = 0
$arr[17,15] = 1
$arr[17,16] = 1
$arr[17,17] = 0
$arr[18,0] = 2554
$arr[18,1] = 25
$arr[18,2] = "Totalisateur Volumes évacués ecluse"
$arr[18,3] = 0
$arr[18,4] = 1
$arr[18,5] = 0
$arr[18,6] = 1
$arr[18,7] = 0
$arr[18,8] = "m³"
$arr[18,9] = "B_BARRAGE_VOLUME"
$arr[18,10] = "C_BARRAGE_VOLUME"
$arr[18,11] = 3
$arr[18,12] = 2
$arr[18,13] = 2
$arr[18,14] = 0
$arr[18,15] = 1
$arr[18,16] = 1
$arr[18,17] = 0
$arr[19,0] = 3000
$arr[19,1] = 8
$arr[19,2] = "Debit"
$arr[19,3] = 0
$arr[19,4] = 1
$arr[19,5] = 0
$arr[19,6] = 1
$arr[19,7] = 2
$arr[19,8] = "m³/s"
$arr[19,9] = "B_LANGON_DEBIT"
$arr[19,10] = "C_LANGON_DEBIT"
$arr[19,11] = 7
$arr[19,12] = 1
$arr[19,13] = 2
$arr[19,14] = 0
$arr[19,15] = 3
$arr[19,16] = 9
$arr[19,17] = 0
$arr[20,0] = 9560
$arr[20,1] = 25
$arr[20,2] = "débit total"
$arr[20,3] = 0
$arr[20,4] = 1
$arr[20,5] = 0
$arr[20,6] = 1
$arr[20,7] = 2
$arr[20,8] = $null
$arr[20,9] = "B_BARRAGE_DEBIT"
$arr[20,10] = "C_BARRAGE_DEBIT"
$arr[20,11] = 16
$arr[20,12] = 0
$arr[20,13] = 0
$arr[20,14] = 0
$arr[20,15] = 0
$arr[20,16] = 0
$arr[20,17] = 0
$arr[21,0] = 30020
$arr[21,1] = 30
$arr[21,2] = "débit total"
$arr[21,3] = 0
$arr[21,4] = 1
$arr[21,5] = 0
$arr[21,6] = 1
$arr[21,7] = 2
$arr[21,8] = $null
$arr[21,9] = "B_ALLAIRE_MESURE"
$arr[21,10] = "C_ALLAIRE_MESURE"
$arr[21,11] = 16
$arr[21,12] = 0
$arr[21,13] = 0
$arr[21,14] = 0
$arr[21,15] = 0
$arr[21,16] = 0
$arr[21,17] = 0

$wsDebit.Range("A1:R22").Value = $arr

# --- "debit" sheet becomes the active tab, with E29 selected ---
$wsDebit.Activate()
$wsDebit.Range("E29").Select()

